$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "bcsstk12.mtx"
$ws.Range("D2").Value = 654926088.2605356
$ws.Range("E2").Value = 28
$ws.Range("F2").Value = 0.02260899543762207
$ws.Range("G2").Value = 1473
$ws.Range("B3").Value = "bcsstk12.mtx"
$ws.Range("D3").Value = 654950577.7596854
$ws.Range("E3").Value = 18
$ws.Range("F3").Value = 0.01457500457763672
$ws.Range("G3").Value = 1473
$ws.Range("B4").Value = "bcsstk10.mtx"
$ws.Range("D4").Value = 44724660.27811407
$ws.Range("E4").Value = 55
$ws.Range("F4").Value = 0.02274203300476074
$ws.Range("G4").Value = 1086
$ws.Range("B5").Value = "bcsstk10.mtx"
$ws.Range("D5").Value = 44735156.50494859
$ws.Range("E5").Value = 44
$ws.Range("F5").Value = 0.01542854309082031
$ws.Range("G5").Value = 1086
$ws.Range("B6").Value = "bcsstk13.mtx"
$ws.Range("D6").Value = 3112885321254.005
$ws.Range("E6").Value = 176
$ws.Range("F6").Value = 0.3588168621063232
$ws.Range("G6").Value = 2003
$ws.Range("B7").Value = "bcsstk13.mtx"
$ws.Range("D7").Value = 3115651354786.666
$ws.Range("E7").Value = 130
$ws.Range("F7").Value = 0.2733302116394043
$ws.Range("G7").Value = 2003
$ws.Range("B8").Value = "bcsstk15.mtx"
$ws.Range("D8").Value = 6537666522.476533
$ws.Range("E8").Value = 79
$ws.Range("F8").Value = 0.6137120723724365
$ws.Range("G8").Value = 3948
$ws.Range("B9").Value = "bcsstk15.mtx"
$ws.Range("D9").Value = 6538214586.124658
$ws.Range("E9").Value = 59
$ws.Range("F9").Value = 0.4607248306274414
$ws.Range("G9").Value = 3948
$ws.Range("B10").Value = "bcsstk07.mtx"
$ws.Range("D10").Value = 3485075325.977184
$ws.Range("E10").Value = 21
$ws.Range("F10").Value = 0.00222468376159668
$ws.Range("G10").Value = 420
$ws.Range("B11").Value = "bcsstk07.mtx"
$ws.Range("D11").Value = 3485103329.846988
$ws.Range("E11").Value = 14
$ws.Range("F11").Value = 0.001317501068115234
$ws.Range("G11").Value = 420
$ws.Range("B12").Value = "bcsstk09.mtx"
$ws.Range("D12").Value = 67579067.21497774
$ws.Range("E12").Value = 220
$ws.Range("F12").Value = 0.07413887977600098
$ws.Range("G12").Value = 1083
$ws.Range("B13").Value = "bcsstk09.mtx"
$ws.Range("D13").Value = 67845967.03713822
$ws.Range("E13").Value = 90
$ws.Range("F13").Value = 0.02515172958374023
$ws.Range("G13").Value = 1083
$ws.Range("B14").Value = "bcsstk11.mtx"
$ws.Range("D14").Value = 654926088.2605356
$ws.Range("E14").Value = 28
$ws.Range("F14").Value = 0.03079676628112793
$ws.Range("G14").Value = 1473
$ws.Range("B15").Value = "bcsstk11.mtx"
$ws.Range("D15").Value = 654950577.7596854
$ws.Range("E15").Value = 18
$ws.Range("F15").Value = 0.01926302909851074
$ws.Range("G15").Value = 1473
$ws.Range("B16").Value = "bcsstk16.mtx"
$ws.Range("D16").Value = 4085202504.101666
$ws.Range("E16").Value = 35
$ws.Range("F16").Value = 0.4120795726776123
$ws.Range("G16").Value = 4884
$ws.Range("B17").Value = "bcsstk16.mtx"
$ws.Range("D17").Value = 4085245570.077303
$ws.Range("E17").Value = 27
$ws.Range("F17").Value = 0.323805570602417
$ws.Range("G17").Value = 4884
$ws.Range("B18").Value = "bcsstk06.mtx"
$ws.Range("D18").Value = 3485075325.977184
$ws.Range("E18").Value = 21
$ws.Range("F18").Value = 0.002191543579101562
$ws.Range("G18").Value = 420
$ws.Range("B19").Value = "bcsstk06.mtx"
$ws.Range("D19").Value = 3485103329.846988
$ws.Range("E19").Value = 14
$ws.Range("F19").Value = 0.00138545036315918
$ws.Range("G19").Value = 420
$ws.Range("B20").Value = "bcsstk08.mtx"
$ws.Range("D20").Value = 76569996762.00502
$ws.Range("E20").Value = 10
$ws.Range("F20").Value = 0.005293369293212891
$ws.Range("G20").Value = 1074
$ws.Range("B21").Value = "bcsstk08.mtx"
$ws.Range("D21").Value = 76570318161.79057
$ws.Range("E21").Value = 9
$ws.Range("F21").Value = 0.00421452522277832
$ws.Range("G21").Value = 1074
